$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '30.258.40'
$dCell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.60%  '

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.864.43'
$dCell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.30%  '

$dCell = $ws.Cells.Item(4, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '236.68'
$dCell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.50%  '

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.000'
$dCell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.00%  '

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.4679'
$dCell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.58%  '

$dCell = $ws.Cells.Item(8, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.2856'
$dCell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  +1.96%  '

$dCell = $ws.Cells.Item(9, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.06534'
$dCell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.18%  '

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '22.27'
$dCell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +14.40%  '

$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.07911'
$dCell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +1.27%  '

$dCell = $ws.Cells.Item(12, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '97.72'
$dCell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +1.51%  '

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.868.36'
$dCell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.50%  '

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.174'
$dCell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.19%  '

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.6832'
$dCell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  +3.01%  '

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '279.43'
$dCell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.14%  '

$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '30.271.37'
$dCell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.52%  '

$dCell = $ws.Cells.Item(18, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '13.68'
$dCell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +9.17%  '

$ws.Cells.Item(19, 5).Value = '  -0.10%  '

$dCell = $ws.Cells.Item(20, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '5.399'
$dCell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -1.42%  '

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.000007332'
$dCell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.75%  '

$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '2.113.38'
$dCell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.55%  '

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.9995'
$dCell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -0.07%  '

$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '6.172'
$dCell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +1.04%  '

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '167.66'
$dCell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.72%  '

$dCell = $ws.Cells.Item(26, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '9.259'
$dCell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -0.23%  '

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '19.14'
$dCell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +1.91%  '

$dCell = $ws.Cells.Item(28, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.939'
$dCell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +2.13%  '

$ws.Cells.Item(29, 5).Value = '  +3.97%  '

$dCell = $ws.Cells.Item(30, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.09839'
$dCell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +3.34%  '

$dCell = $ws.Cells.Item(31, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '4.393'
$dCell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -0.48%  '

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.482'
$dCell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +1.35%  '

$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '4.065'
$dCell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.36%  '

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.04744'
$dCell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.50%  '

$ws.Cells.Item(35, 5).Value = '  +4.40%  '

$dCell = $ws.Cells.Item(36, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.7116'
$dCell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  +1.88%  '

$ws.Cells.Item(37, 5).Value = '  +0.16%  '

$ws.Cells.Item(38, 5).Value = '  +1.94%  '

$ws.Cells.Item(39, 5).Value = '  +4.36%  '

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '76.73'
$dCell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +5.66%  '

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '6.295'
$dCell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +0.29%  '

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '1.960'
$dCell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +2.69%  '

$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.8512'
$dCell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -0.20%  '

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.4182'
$dCell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +1.21%  '

$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.9992'
$dCell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -0.05%  '

$dCell = $ws.Cells.Item(46, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '103.49'
$dCell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -0.37%  '

$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '962.09'
$dCell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.13%  '

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '7.227'
$dCell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.92%  '

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '9.345'
$dCell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.43%  '

$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '34.22'
$dCell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.91%  '

$dCell = $ws.Cells.Item(51, 4)
$dCell.NumberFormat = "@"
$dCell.Value = '0.05647'
$dCell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  +0.51%  '
